$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the style from an existing header cell (G1) so H1 reuses the same
# cellXf (bold font, centered, thin border) instead of minting a new style.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Fill in the new "Save" data column (H2:H7) per the diff.
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 1
$ws.Range("H4").Value = 1
$ws.Range("H5").Value = 1
$ws.Range("H6").Value = 0
$ws.Range("H7").Value = 0
